$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.7231684923172
$ws.Range("B1").Value = 2.682357311248779
$ws.Range("C1").Value = 2.873579740524292
$ws.Range("D1").Value = 3.239982128143311
$ws.Range("E1").Value = 2.600570201873779
